$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "experiment_description": add experiment 3 (base model) and
# experiment 4 (extended model), both reusing the "Test experiment" name.
# ---------------------------------------------------------------------------
$wsDesc = $wb.Worksheets.Item("experiment_description")
$wsDesc.Range("A4").Value = 3
$wsDesc.Range("B4").Value = "Test experiment"
$wsDesc.Range("C4").Value = "base"

$wsDesc.Range("A5").Value = 4
$wsDesc.Range("B5").Value = "Test experiment"
$wsDesc.Range("C5").Value = "extended"

# ---------------------------------------------------------------------------
# Sheet "experiment_specification": the max-splitting-variable change.
# Experiment 2's inpatient_ward / intensive_care_unit transitions are split
# into green/red variants (like home already was), and two new experiments
# (3 and 4) are added: 3 mirrors experiment 1 but splits the max variable
# for home into "age_three"; 4 mirrors experiment 2 (green/red) but also
# uses "age_three" for the home transitions.
# ---------------------------------------------------------------------------
$wsSpec = $wb.Worksheets.Item("experiment_specification")

# Replace rows 16-19 (old, single inpatient_ward/icu rows for experiment 2)
# with the new 16-21 (green/red split) block.
$wsSpec.Range("A16").Value = 2
$wsSpec.Range("B16").Value = "transition"
$wsSpec.Range("C16").Value = "inpatient_ward-green"
$wsSpec.Range("D16").Value = "none"

$wsSpec.Range("A17").Value = 2
$wsSpec.Range("B17").Value = "transition"
$wsSpec.Range("C17").Value = "inpatient_ward-red"
$wsSpec.Range("D17").Value = "none"

$wsSpec.Range("A18").Value = 2
$wsSpec.Range("B18").Value = "length_of_stay"
$wsSpec.Range("C18").Value = "inpatient_ward"
$wsSpec.Range("D18").Value = "none"

$wsSpec.Range("A19").Value = 2
$wsSpec.Range("B19").Value = "transition"
$wsSpec.Range("C19").Value = "intensive_care_unit-green"
$wsSpec.Range("D19").Value = "none"

$wsSpec.Range("A20").Value = 2
$wsSpec.Range("B20").Value = "transition"
$wsSpec.Range("C20").Value = "intensive_care_unit-red"
$wsSpec.Range("D20").Value = "none"

$wsSpec.Range("A21").Value = 2
$wsSpec.Range("B21").Value = "length_of_stay"
$wsSpec.Range("C21").Value = "intensive_care_unit"
$wsSpec.Range("D21").Value = "none"

# Experiment 3: same layout as experiment 1, but the home transition now
# uses "age_three" as the splitting variable.
$wsSpec.Range("A22").Value = 3
$wsSpec.Range("B22").Value = "heuristic"
$wsSpec.Range("D22").Value = "heuristic_1"

$wsSpec.Range("A23").Value = 3
$wsSpec.Range("B23").Value = "heuristic"
$wsSpec.Range("D23").Value = "heuristic_2"

$wsSpec.Range("A24").Value = 3
$wsSpec.Range("B24").Value = "heuristic"
$wsSpec.Range("D24").Value = "heuristic_3"

$wsSpec.Range("A25").Value = 3
$wsSpec.Range("B25").Value = "transition"
$wsSpec.Range("C25").Value = "home"
$wsSpec.Range("D25").Value = "age_three"

$wsSpec.Range("A26").Value = 3
$wsSpec.Range("B26").Value = "length_of_stay"
$wsSpec.Range("C26").Value = "home"
$wsSpec.Range("D26").Value = "age_simple"

$wsSpec.Range("A27").Value = 3
$wsSpec.Range("B27").Value = "transition"
$wsSpec.Range("C27").Value = "inpatient_ward"
$wsSpec.Range("D27").Value = "age_simple"

$wsSpec.Range("A28").Value = 3
$wsSpec.Range("B28").Value = "length_of_stay"
$wsSpec.Range("C28").Value = "inpatient_ward"
$wsSpec.Range("D28").Value = "none"

$wsSpec.Range("A29").Value = 3
$wsSpec.Range("B29").Value = "transition"
$wsSpec.Range("C29").Value = "intensive_care_unit"
$wsSpec.Range("D29").Value = "age_simple"

$wsSpec.Range("A30").Value = 3
$wsSpec.Range("B30").Value = "length_of_stay"
$wsSpec.Range("C30").Value = "intensive_care_unit"
$wsSpec.Range("D30").Value = "none"

# Experiment 4: same layout as experiment 2 (green/red split), but the
# home-green/home-red transitions use "age_three" as the splitting variable.
$wsSpec.Range("A31").Value = 4
$wsSpec.Range("B31").Value = "heuristic"
$wsSpec.Range("D31").Value = "heuristic_1"

$wsSpec.Range("A32").Value = 4
$wsSpec.Range("B32").Value = "transition"
$wsSpec.Range("C32").Value = "home-green"
$wsSpec.Range("D32").Value = "age_three"

$wsSpec.Range("A33").Value = 4
$wsSpec.Range("B33").Value = "transition"
$wsSpec.Range("C33").Value = "home-red"
$wsSpec.Range("D33").Value = "age_three"

$wsSpec.Range("A34").Value = 4
$wsSpec.Range("B34").Value = "length_of_stay"
$wsSpec.Range("C34").Value = "home-green"
$wsSpec.Range("D34").Value = "age_simple"

$wsSpec.Range("A35").Value = 4
$wsSpec.Range("B35").Value = "length_of_stay"
$wsSpec.Range("C35").Value = "home-red"
$wsSpec.Range("D35").Value = "age_simple"

$wsSpec.Range("A36").Value = 4
$wsSpec.Range("B36").Value = "transition"
$wsSpec.Range("C36").Value = "inpatient_ward-green"
$wsSpec.Range("D36").Value = "none"

$wsSpec.Range("A37").Value = 4
$wsSpec.Range("B37").Value = "transition"
$wsSpec.Range("C37").Value = "inpatient_ward-red"
$wsSpec.Range("D37").Value = "none"

$wsSpec.Range("A38").Value = 4
$wsSpec.Range("B38").Value = "length_of_stay"
$wsSpec.Range("C38").Value = "inpatient_ward"
$wsSpec.Range("D38").Value = "none"

$wsSpec.Range("A39").Value = 4
$wsSpec.Range("B39").Value = "transition"
$wsSpec.Range("C39").Value = "intensive_care_unit-green"
$wsSpec.Range("D39").Value = "none"

$wsSpec.Range("A40").Value = 4
$wsSpec.Range("B40").Value = "transition"
$wsSpec.Range("C40").Value = "intensive_care_unit-red"
$wsSpec.Range("D40").Value = "none"

$wsSpec.Range("A41").Value = 4
$wsSpec.Range("B41").Value = "length_of_stay"
$wsSpec.Range("C41").Value = "intensive_care_unit"
$wsSpec.Range("D41").Value = "none"

# ---------------------------------------------------------------------------
# Sheet "run_specification": add two more heuristic rows for run 2.
# ---------------------------------------------------------------------------
$wsRunSpec = $wb.Worksheets.Item("run_specification")
$wsRunSpec.Range("A5").Value = 2
$wsRunSpec.Range("B5").Value = 3

$wsRunSpec.Range("A6").Value = 2
$wsRunSpec.Range("B6").Value = 4

# ---------------------------------------------------------------------------
# Selections / active sheet, matching the saved workbook view.
# ---------------------------------------------------------------------------
$wsDesc.Range("C5").Select()
$wsSpec.Range("A31:A41").Select()
$wsRunSpec.Range("B6").Select()

$wsRunSpec.Activate()
